$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (23) describing the abs() function, following the pattern
# used by existing rows (함수명/동작/설명/예시 columns).
$ws.Range("A23").Value = "abs()"
$ws.Range("B23").Value = "주어진 숫자의 절댓값을 반환"
$ws.Range("C23").Value = "정수나 실수: 주어진 숫자가 양수이면 그대로, 음수이면 양수로 반환합니다.`n복소수: 복소수의 크기(절댓값)를 반환합니다."
$ws.Range("D23").Value = "abs_value = abs(x)"

# Match the styling used by similar rows (row 17 uses the same column layout:
# style 1 for A/B/D, style 2 -- wrap text -- for C).
$ws.Range("A23").Style = $ws.Range("A17").Style
$ws.Range("B23").Style = $ws.Range("B17").Style
$ws.Range("C23").Style = $ws.Range("C17").Style
$ws.Range("D23").Style = $ws.Range("D17").Style

$ws.Rows.Item(23).RowHeight = 40

# Update the selection / view state to reflect what was left selected after
# the edit.
$ws.Range("C2").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
